$d = $word.ActiveDocument

# Locate the paragraph "Ilvl of items on inspect panel" - the new entries
# are inserted as new paragraphs right after it.
$anchor = $null
foreach ($p in $d.Paragraphs) {
    $raw = $p.Range.Text
    $trimmed = $raw.TrimEnd([char]13, [char]7)
    if ($trimmed -eq "Ilvl of items on inspect panel") {
        $anchor = $p
        break
    }
}

if ($anchor -ne $null) {
    $newLines = @(
        "Option to tenable/disable BossesKilled",
        "Droprate on items tooltip",
        "New option panel for secondaries options"
    )

    $cur = $anchor
    foreach ($line in $newLines) {
        $cur.Range.InsertParagraphAfter()
        $cur = $cur.Next()
        $cur.Range.InsertAfter($line)
    }
}
